$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 64780
$ws.Range("J57").Value = 64780
$ws.Range("L57").Value = 194340
$ws.Range("N57").Value = -195338

$ws.Range("H76").Value = 35486696
$ws.Range("I76").Value = 39288520
$ws.Range("K76").Value = 39288520
$ws.Range("M76").Value = -39288205

$ws.Range("H79").Value = 35486696
$ws.Range("I79").Value = 39288520
$ws.Range("K79").Value = 39288520
$ws.Range("M79").Value = -39287428

$ws.Range("H86").Value = 111115736
$ws.Range("I86").Value = 3608
$ws.Range("K86").Value = 3608
$ws.Range("M86").Value = -2485

$ws.Range("H89").Value = 111115736
$ws.Range("I89").Value = 3608
$ws.Range("K89").Value = 18040
$ws.Range("M89").Value = -12424

$ws.Range("H92").Value = 1211.2858
$ws.Range("I92").Value = 1303.7894
$ws.Range("K92").Value = 1303.7894
$ws.Range("M92").Value = -55.78939999999989

$ws.Range("H117").Value = 17685.5
$ws.Range("J117").Value = 17685.5
$ws.Range("L117").Value = 17685.5
$ws.Range("N117").Value = -26863.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2283.2083
$ws.Range("I2").Value = 1906.9333
$ws.Range("K2").Value = 1906.9333
$ws.Range("M2").Value = -1793.9333

$ws.Range("H32").Value = 11240.507
$ws.Range("I32").Value = 6516.846
$ws.Range("K32").Value = 6516.846
$ws.Range("M32").Value = -6229.846

$ws.Range("H45").Value = 1140
$ws.Range("I45").Value = 933.3333
$ws.Range("J45").Value = 1450
$ws.Range("K45").Value = 933.3333
$ws.Range("L45").Value = 1450
$ws.Range("M45").Value = -556.3333
$ws.Range("N45").Value = -2204

$ws.Range("H61").Value = 2807.8572
$ws.Range("I61").Value = 2658.2
$ws.Range("K61").Value = 2658.2
$ws.Range("M61").Value = -2446.2

$ws.Range("H74").Value = 1582.7354
$ws.Range("I74").Value = 1380.4482
$ws.Range("J74").Value = 2756
$ws.Range("K74").Value = 1380.4482
$ws.Range("L74").Value = 2756
$ws.Range("M74").Value = -506.4482
$ws.Range("N74").Value = -4504

$ws.Range("H77").Value = 1582.7354
$ws.Range("I77").Value = 1380.4482
$ws.Range("J77").Value = 2756
$ws.Range("K77").Value = 6902.241
$ws.Range("L77").Value = 13780
$ws.Range("M77").Value = -2534.241
$ws.Range("N77").Value = -22516

$ws.Range("H101").Value = 40000
$ws.Range("J101").Value = 40000
$ws.Range("L101").Value = 40000
$ws.Range("N101").Value = -46490

$ws.Range("H116").Value = 2283.2083
$ws.Range("I116").Value = 1906.9333
$ws.Range("K116").Value = 1906.9333
$ws.Range("M116").Value = 387.0667000000001

$ws.Range("H136").Value = 2807.8572
$ws.Range("I136").Value = 2658.2
$ws.Range("K136").Value = 7974.599999999999
$ws.Range("M136").Value = -5424.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2283.2083
$ws.Range("I3").Value = 1906.9333
$ws.Range("K3").Value = 1906.9333
$ws.Range("M3").Value = -1792.9333

$ws.Range("H57").Value = 75935
$ws.Range("J57").Value = 75935
$ws.Range("L57").Value = 75935
$ws.Range("N57").Value = -77375

$ws.Range("H86").Value = 55558600
$ws.Range("I86").Value = 76925370
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 76925370
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -76924247
$ws.Range("N86").Value = -7246

$ws.Range("H89").Value = 55558600
$ws.Range("I89").Value = 76925370
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 384626850
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -384621234
$ws.Range("N89").Value = -36232

$ws.Range("H117").Value = 48000
$ws.Range("J117").Value = 48000
$ws.Range("L117").Value = 48000
$ws.Range("N117").Value = -57178

$ws.Range("H136").Value = 75935
$ws.Range("J136").Value = 75935
$ws.Range("L136").Value = 75935
$ws.Range("N136").Value = -86135

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 69095
$ws.Range("J138").Value = 69095
$ws.Range("L138").Value = 69095
$ws.Range("N138").Value = -79375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1070.1904
$ws.Range("J122").Value = 1137.4445
$ws.Range("L122").Value = 10237.0005
$ws.Range("N122").Value = -15137.0005

$ws.Range("H137").Value = 6601.3687
$ws.Range("I137").Value = 3438.3333
$ws.Range("J137").Value = 8061.231
$ws.Range("K137").Value = 10314.9999
$ws.Range("L137").Value = 24183.693
$ws.Range("M137").Value = -5214.999899999999
$ws.Range("N137").Value = -34383.693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5466.6665
$ws.Range("I70").Value = 5333.3335
$ws.Range("K70").Value = 5333.3335
$ws.Range("M70").Value = -5063.3335

$ws.Range("H73").Value = 5466.6665
$ws.Range("I73").Value = 5333.3335
$ws.Range("K73").Value = 5333.3335
$ws.Range("M73").Value = -4397.3335

$ws.Range("H126").Value = 5414.75
$ws.Range("I126").Value = 5903.2
$ws.Range("J126").Value = 4600.6665
$ws.Range("K126").Value = 17709.6
$ws.Range("L126").Value = 13801.9995
$ws.Range("M126").Value = -15239.6
$ws.Range("N126").Value = -18741.9995

$ws.Range("H132").Value = 3192.0833
$ws.Range("I132").Value = 3022.889
$ws.Range("J132").Value = 3699.6667
$ws.Range("K132").Value = 9068.667000000001
$ws.Range("L132").Value = 11099.0001
$ws.Range("M132").Value = -6538.667000000001
$ws.Range("N132").Value = -16159.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1889.8286
$ws.Range("I7").Value = 1928.3334
$ws.Range("J7").Value = 1805.8182
$ws.Range("K7").Value = 1928.3334
$ws.Range("L7").Value = 1805.8182
$ws.Range("M7").Value = -1816.3334
$ws.Range("N7").Value = -2029.8182

$ws.Range("H126").Value = 1889.8286
$ws.Range("I126").Value = 1928.3334
$ws.Range("J126").Value = 1805.8182
$ws.Range("K126").Value = 5785.0002
$ws.Range("L126").Value = 5417.4546
$ws.Range("M126").Value = -3315.0002
$ws.Range("N126").Value = -10357.4546

$ws.Range("H136").Value = 3935.6562
$ws.Range("I136").Value = 1851.7084
$ws.Range("K136").Value = 5555.1252
$ws.Range("M136").Value = -3005.1252

$ws.Range("H139").Value = 59682.5
$ws.Range("J139").Value = 58715
$ws.Range("L139").Value = 58715
$ws.Range("N139").Value = -68995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 6322
$ws.Range("J45").Value = 6646.4
$ws.Range("L45").Value = 6646.4
$ws.Range("N45").Value = -7628.4

$ws.Range("H105").Value = 40333.332
$ws.Range("J105").Value = 40333.332
$ws.Range("L105").Value = 40333.332
$ws.Range("N105").Value = -47321.332
